# Updates the cryptocurrency price/volume table to the latest scrape.
# Mirrors a plain re-assignment of the changed <is><t> cell contents;
# price strings that would otherwise be auto-detected as numbers are
# entered with a leading apostrophe (quote-prefix) so Excel keeps them
# as literal text, exactly like the source inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.986.91"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "3.536.53"
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'" + "188.60"
$ws.Range("E5").Value = "  +7.97%  "
$ws.Range("D6").Value = "'" + "562.94"
$ws.Range("E6").Value = "  +6.52%  "
$ws.Range("D7").Value = "'" + "0.626"
$ws.Range("E7").Value = "  +4.94%  "
$ws.Range("D8").Value = "3.528.43"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("E11").Value = "  +15.38%  "
$ws.Range("D12").Value = "'" + "54.69"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "'" + "0.0000276"
$ws.Range("E13").Value = "  +6.94%  "
$ws.Range("D14").Value = "'" + "9.37"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").Value = "4.098.59"
$ws.Range("E15").Value = "  +5.26%  "
$ws.Range("D16").Value = "3.540.81"
$ws.Range("E16").Value = "  +5.56%  "
$ws.Range("D17").Value = "'" + "18.64"
$ws.Range("E17").Value = "  +6.03%  "
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("D19").Value = "67.061.78"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("E20").Value = "  +7.64%  "
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "'" + "421.88"
$ws.Range("E22").Value = "  +12.73%  "
$ws.Range("D23").Value = "'" + "4.12"
$ws.Range("E23").Value = "  +10.62%  "
$ws.Range("D24").Value = "'" + "85.57"
$ws.Range("E24").Value = "  +5.07%  "
$ws.Range("D25").Value = "'" + "4.19"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").Value = "'" + "11.06"
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("D27").Value = "'" + "2.91"
$ws.Range("E27").Value = "  +7.78%  "
$ws.Range("D28").Value = "'" + "12.32"
$ws.Range("E28").Value = "  +8.72%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  +10.16%  "
$ws.Range("D31").Value = "'" + "30.44"
$ws.Range("E31").Value = "  +5.25%  "
$ws.Range("D32").Value = "'" + "636.01"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'" + "6.66"
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").Value = "'" + "11.74"
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("D36").Value = "'" + "60.44"
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("D37").Value = "0.0₃0826"
$ws.Range("E37").Value = "  +12.13%  "
$ws.Range("E38").Value = "  +19.23%  "
$ws.Range("D39").Value = "'" + "38.39"
$ws.Range("E39").Value = "  +5.25%  "
$ws.Range("D40").Value = "'" + "0.998"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'" + "3.36"
$ws.Range("E42").Value = "  +12.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.138.78"
$ws.Range("E43").Value = "  +5.57%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  +9.78%  "
$ws.Range("E47").Value = "  +10.26%  "
$ws.Range("D48").Value = "'" + "0.0418"
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").Value = "'" + "0.133"
$ws.Range("E50").Value = "  +5.71%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'" + "140.26"
$ws.Range("E51").Value = "  +2.02%  "
